$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset "Iwantoro et al. (2019)" (rows 3-5, Numerical / M2) was removed
# from the "Overall results" sheet. Deleting the rows shifts everything below
# up by three and Excel/the engine automatically:
#   - renumbers row references
#   - shrinks the used-range dimension (A1:H40 -> A1:H37)
#   - adjusts the merged blank separator row (A39:H39 -> A36:H36)
#   - drops the now-unused "Iwantoro et al. (2019)" shared string
$ws.Rows("3:5").Delete() | Out-Null

# Reflect where the author's cursor ended up after the edit.
$ws.Range("D20").Select() | Out-Null
